{"js": "// Locate the target paragraph (the comment line that gets extended) and the\n// new lines of C++ code/comments that were added after it, just before the\n// trailing empty paragraph at the very end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText =\n  \"                    // Aqu\u00ed puedes implementar la l\u00f3gica para convertir el tiempo en segundos a\";\nconst appendedSuffix = \" un formato de hora\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === oldText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph with text: \" + oldText);\n}\n\n// Extend the existing comment line with the missing words.\nanchor.insertText(appendedSuffix, Word.InsertLocation.end);\n\n// Lines to insert, in order, right after the anchor paragraph. `null`\n// represents a blank paragraph (no text run at all).\nconst newLines = [\n  \"                    // y sumarle la hora de inicio proporcionada.\",\n  '                    cout << \"El tiempo estimado de llegada es: \" << tiempoLlegada << \" segundos\\\\n\";',\n  \"                }\",\n  \"                break;\",\n  \"            }\",\n  \"            case 0: {\",\n  '                cout << \"Saliendo del simulador.\\\\n\";',\n  \"                break;\",\n  \"            }\",\n  \"            default: {\",\n  '                cout << \"Opci\u00f3n no v\u00e1lida. Por favor, ingrese una opci\u00f3n v\u00e1lida.\\\\n\";',\n  \"                break;\",\n  \"            }\",\n  \"        }\",\n  \"    } while (opcion != 0);\",\n  null,\n  \"    return 0;\",\n  \"}\",\n];\n\nlet insertAfter = anchor;\nfor (const line of newLines) {\n  const inserted = insertAfter.insertParagraph(line === null ? \"\" : line, Word.InsertLocation.after);\n  insertAfter = inserted;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the anchor paragraph - the comment line that gets extended with\n# \"un formato de hora\".\n$oldText = \"                    // Aqu\u00ed puedes implementar la l\u00f3gica para convertir el tiempo en segundos a\"\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    # Range.Text includes the trailing paragraph mark (CR); strip it for comparison.\n    $t = $t.TrimEnd([char]13)\n    if ($t -eq $oldText) {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find anchor paragraph\"\n}\n\n# Extend the existing comment line with the missing words.\n$anchor.Range.InsertAfter(\" un formato de hora\")\n\n# Lines to insert, in order, right after the anchor paragraph. $null\n# represents a blank paragraph (no text).\n$newLines = @(\n    \"                    // y sumarle la hora de inicio proporcionada.\",\n    '                    cout << \"El tiempo estimado de llegada es: \" << tiempoLlegada << \" segundos\\n\";',\n    \"                }\",\n    \"                break;\",\n    \"            }\",\n    \"            case 0: {\",\n    '                cout << \"Saliendo del simulador.\\n\";',\n    \"                break;\",\n    \"            }\",\n    \"            default: {\",\n    '                cout << \"Opci\u00f3n no v\u00e1lida. Por favor, ingrese una opci\u00f3n v\u00e1lida.\\n\";',\n    \"                break;\",\n    \"            }\",\n    \"        }\",\n    \"    } while (opcion != 0);\",\n    $null,\n    \"    return 0;\",\n    \"}\"\n)\n\n$insertRange = $anchor.Range\nforeach ($line in $newLines) {\n    $insertRange.InsertParagraphAfter()\n    $insertRange = $insertRange.Next(4, 1)\n    if ($line -ne $null) {\n        $insertRange.Text = $line\n    }\n}\n"}
